$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.729.54'
$ws.Range("E2").Value = '  +2.26%  '

# Row 3
$ws.Range("D3").Value = '1.816.76'
$ws.Range("E3").Value = '  -0.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.10'
$ws.Range("E5").Value = '  -2.71%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.14%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4436'
$ws.Range("E7").Value = '  +4.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3821'
$ws.Range("E8").Value = '  +8.85%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.82'
$ws.Range("E9").Value = '  -1.82%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07757'
$ws.Range("E10").Value = '  +3.95%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.156'
$ws.Range("E11").Value = '  +0.40%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.92'
$ws.Range("E12").Value = '  -0.47%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9995'
$ws.Range("E13").Value = '  -0.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.368'
$ws.Range("E14").Value = '  +1.47%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.618'
$ws.Range("E15").Value = '  +4.47%  '

# Row 16
$ws.Range("D16").Value = '1.812.51'
$ws.Range("E16").Value = '  +0.05%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001100'
$ws.Range("E17").Value = '  +1.19%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06752'
$ws.Range("E18").Value = '  +0.90%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.81'
$ws.Range("E19").Value = '  -0.33%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9982'
$ws.Range("E20").Value = '  -0.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.89'
$ws.Range("E21").Value = '  +3.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.365'
$ws.Range("E22").Value = '  -1.49%  '

# Row 23
$ws.Range("D23").Value = '28.712.82'
$ws.Range("E23").Value = '  +2.19%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.95'
$ws.Range("E24").Value = '  -0.76%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.437'
$ws.Range("E25").Value = '  +1.98%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.84'
$ws.Range("E26").Value = '  +0.37%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.405'
$ws.Range("E27").Value = '  -3.75%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.45'
$ws.Range("E28").Value = '  -1.82%  '

# Row 29
$ws.Range("D29").Value = '2.018.14'
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.295'
$ws.Range("E30").Value = '  -0.78%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.38'
$ws.Range("E31").Value = '  +0.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.972'
$ws.Range("E32").Value = '  -2.26%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.909'
$ws.Range("E33").Value = '  -0.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09337'
$ws.Range("E34").Value = '  +1.58%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2287'
$ws.Range("E35").Value = '  +5.32%  '

# Row 36
$ws.Range("E36").Value = '  +0.12%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06422'
$ws.Range("E37").Value = '  +2.28%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02362'
$ws.Range("E38").Value = '  -0.12%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.285'
$ws.Range("E39").Value = '  +0.52%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6732'
$ws.Range("E40").Value = '  -0.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.214'
$ws.Range("E41").Value = '  -0.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.242'
$ws.Range("E42").Value = '  +1.41%  '

# Row 43
$ws.Range("E43").Value = '  -2.90%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.11'
$ws.Range("E44").Value = '  +0.18%  '

# Row 45
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9978'
$ws.Range("E45").Value = '  -0.20%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6162'
$ws.Range("E46").Value = '  +0.35%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.813'
$ws.Range("E47").Value = '  -1.68%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.71'
$ws.Range("E48").Value = '  +0.86%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.069'
$ws.Range("E49").Value = '  +0.94%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07132'
$ws.Range("E50").Value = '  +0.36%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.165'
$ws.Range("E51").Value = '  -1.33%  '
